# Apply "Nodata (-1)" fill-in to all Body rows that were missing the
# Equatorial/Polar radii columns (C:F), matching values already present
# for rows that have real CSV-sourced radii data.
#
# Each block below is a contiguous run of worksheet rows (grouped from the
# diff) that previously had only columns A (Naif_id) and B (Body) filled
# in. We set C:F to -1 and format them with the same "integer" number
# format Excel assigns to whole-number cells (numFmtId 1 == "0"), which
# creates the third cellXfs entry alongside the workbook's existing
# General (0) and 2-decimal (1) styles.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$noDataRanges = @(
    "C26:F59",
    "C79:F91",
    "C95:F113",
    "C130:F141",
    "C151:F155",
    "C158:F161",
    "C164:F164",
    "C166:F167",
    "C171:F178"
)

foreach ($addr in $noDataRanges) {
    $rng = $ws.Range($addr)
    $rng.Value = -1
    $rng.NumberFormat = "0"
}

# Row 169 already has D/E/F (real data from the CSV); only C was missing.
$c169 = $ws.Range("C169")
$c169.Value = -1
$c169.NumberFormat = "0"

# Reflect the author's on-screen state when the file was saved: zoomed to
# 85% with C171:F178 (the last block of newly-filled Nodata rows)
# selected/active.
$excel.ActiveWindow.Zoom = 85
$ws.Range("C171:F178").Select()
